$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $xml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    [void]$r.InsertXML($xml)
}

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>ContosoLearn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>竞争对手</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> SWOT</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 1 $xml1

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Fabrikam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> Learning</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>：</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 2 $xml2

$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>优势：</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Fabrikam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> Learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>提供了一套全面的分析和报告工具。它可确保持续监控教学和学习活动，并查明需要解决的问题区域。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 3 $xml3

$xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>缺点：</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>虽然</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Fabrikam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> Learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>具有强大的报告功能，但由于其全面性，可能会让一些用户不知所措。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 4 $xml4

$xml5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>机会：</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>对个性化学习体验和数据驱动建议的需求越来越大。</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Fabrikam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> Learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>可以利用其强大的分析和报告工具来满足此需求。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 5 $xml5

$xml6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>威胁：</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">eLearning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>市场的激烈竞争，许多公司都提供类似的功能。</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Fabrikam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> Learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>需要持续创新才能保持领先。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 6 $xml6

$xml7 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>AdatumLearn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>：</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 7 $xml7

$xml8 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>优势：</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>AdatumLearn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>提供有关</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> MOST </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>和</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> SWOT </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>等业务分析技术的课程。这表明他们致力于向用户提供有价值的内容。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 8 $xml8

$xml9 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>缺点：</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>课程中提供的信息是对第三方所生成信息进行编译后的结果。这可能不如原始内容那么有价值。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 9 $xml9

$xml10 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>机会：</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>AdatumLearn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>可以创建更多原始内容，为用户提供独特的价值。他们还可以扩展课程产品</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>服务，以涵盖更多主题。</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 10 $xml10

$xml11 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:b/><w:bCs/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>威胁：</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>像</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Fabrikam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> Learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>一样，</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>AdatumLearn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>也面临着</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> eLearning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>市场的激烈竞争。他们需要不断改进自己的产品</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>服务，以保持竞争力。</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos (Body)" w:eastAsia="Microsoft YaHei UI" w:hAnsi="Aptos (Body)" w:cs="Aptos (Body)"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>”</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphXml 11 $xml11

# Trailing empty paragraph: set paragraph-mark formatting directly
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$rLast.Font.NameAscii = "Aptos (Body)"
$rLast.Font.NameFarEast = "Microsoft YaHei UI"
$rLast.Font.NameOther = "Aptos (Body)"
$rLast.Font.NameBi = "Aptos (Body)"
$rLast.LanguageIDFarEast = "zh-CN"
